# Insert 3 new data rows before existing row 492, shifting the old
# rows 492..565 down to 495..568, then populate the 3 new rows with
# their data (columns A..T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 492-565 down by three rows.
$ws.Rows("492:494").Insert()

# Common (constant) values shared by every data row in this sheet.
$colA = 5
$colB = "Macroferia Regional de Talca"
$colC = "Maule"
$colE = 7
$colF = "Fruta"
$colG = 100109
$colH = "Uva"
$colI = 100109001
$colJ = "Uva"

function Set-DataRow {
    param(
        [int]$RowNum,
        [double]$D,
        [string]$K,
        [string]$L,
        [double]$M,
        [double]$N,
        [double]$O,
        [double]$P,
        [string]$Q,
        [string]$R,
        [double]$S,
        [double]$T
    )

    $ws.Cells.Item($RowNum, 1).Value = $colA
    $ws.Cells.Item($RowNum, 2).Value = $colB
    $ws.Cells.Item($RowNum, 3).Value = $colC
    $ws.Cells.Item($RowNum, 4).Value = $D
    $ws.Cells.Item($RowNum, 5).Value = $colE
    $ws.Cells.Item($RowNum, 6).Value = $colF
    $ws.Cells.Item($RowNum, 7).Value = $colG
    $ws.Cells.Item($RowNum, 8).Value = $colH
    $ws.Cells.Item($RowNum, 9).Value = $colI
    $ws.Cells.Item($RowNum, 10).Value = $colJ
    $ws.Cells.Item($RowNum, 11).Value = $K
    $ws.Cells.Item($RowNum, 12).Value = $L
    $ws.Cells.Item($RowNum, 13).Value = $M
    $ws.Cells.Item($RowNum, 14).Value = $N
    $ws.Cells.Item($RowNum, 15).Value = $O
    $ws.Cells.Item($RowNum, 16).Value = $P
    $ws.Cells.Item($RowNum, 17).Value = $Q
    $ws.Cells.Item($RowNum, 18).Value = $R
    $ws.Cells.Item($RowNum, 19).Value = $S
    $ws.Cells.Item($RowNum, 20).Value = $T
}

Set-DataRow 492 44946 "Flame Seedless" "Primera" 200 8000 8000 8000 "`$/bandeja 10 kilos" "Provincia de Limarí" 800 10
Set-DataRow 493 44946 "Superior Seedless" "Primera" 220 8000 8000 8000 "`$/bandeja 10 kilos" "Provincia de Limarí" 800 10
Set-DataRow 494 44946 "Superior Seedless" "Primera" 200 13000 13000 13000 "`$/bandeja 18 kilos" "Provincia de San Felipe de Aconcagua" 722 18
